# LARA Breathing Blocks - Errata update
# 1. Add a new "Percussion 2" row to the end of the first (errata) table.
# 2. Remove the stray <w:lastRenderedPageBreak/> that precedes the
#    "Errata from Version 1 (page 1-19)" heading.

$d = $word.ActiveDocument

# --- 1. Append the new errata row -----------------------------------------
$table = $d.Tables(1)
$newRow = $table.Rows.Add()

$newRow.Cells(1).Range.Text = "Percussion 2"
$newRow.Cells(2).Range.Text = "148"
$newRow.Cells(3).Range.Text = "Is the G in beat 3 meant to be a quaver? "
# Cells(4) (the "Answer" column) is left blank, matching the other rows.

# --- 2. Drop the stale lastRenderedPageBreak marker ------------------------
# Find & Replace the heading text in place; Word regenerates the run
# without carrying over the old lastRenderedPageBreak bookkeeping element.
$d.Content.Find.Execute("Errata from Version 1 (page 1-19)", $true, $false,
                         $false, $false, $false, $true, 1, $false,
                         "Errata from Version 1 (page 1-19)", 2)
